$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date value that was updated from
# 2023-09-03 (serial 45172) to 2023-09-06 (serial 45175) for every
# data row (rows 2 through 51).
for ($r = 2; $r -le 51; $r++) {
    $ws.Cells.Item($r, 3).Value = 45175
}
